# The three data rows 14, 15, 16 need their content cyclically rotated:
#   new row14 <- old row16
#   new row15 <- old row14
#   new row16 <- old row15
# (Row numbers / A column "Id" stay put; every other field moves with its record.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @(
    "A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T",
    "U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK",
    "AL","AM","AN","AO","AP","AQ","AR","AS","AT","AU","AV","AW","AX","AY"
)

# Columns whose values must stay plain text even though they may look like
# numbers/dates (e.g. "2023-07-28", or "5") which Excel would otherwise
# silently reinterpret as a date serial / number when assigned.
$textCols = @("Y", "AA", "I")

# Read current contents of the three rows first, so the writes below don't
# clobber values we still need to read.
$row14 = @{}
$row15 = @{}
$row16 = @{}
foreach ($col in $cols) {
    $row14[$col] = $ws.Range("$col`14").Value()
    $row15[$col] = $ws.Range("$col`15").Value()
    $row16[$col] = $ws.Range("$col`16").Value()
}

function Set-CellValue($col, $row, $value) {
    $cell = $ws.Range("$col$row")
    if (($textCols -contains $col) -and ($value -ne $null) -and ("$value" -ne "")) {
        # Force text interpretation, then strip the style bit that
        # NumberFormat leaves behind so the cell stays plain/unstyled.
        $cell.NumberFormat = "@"
        $cell.Value = $value
        $cell.ClearFormats()
    } else {
        $cell.Value = $value
    }
}

foreach ($col in $cols) {
    if ($row14[$col] -ne $row16[$col]) {
        Set-CellValue $col 14 $row16[$col]
    }
    if ($row15[$col] -ne $row14[$col]) {
        Set-CellValue $col 15 $row14[$col]
    }
    if ($row16[$col] -ne $row15[$col]) {
        Set-CellValue $col 16 $row15[$col]
    }
}
